$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D2").Value = "2016-01-15 10:09:48"
$wsZh.Range("G2").Value = "2016-01-15 10:10:32"

$wsDe.Range("D2").Value = "2016-01-15 10:09:58"
$wsDe.Range("G2").Value = "2016-01-15 10:10:50"
